# Apply the "Natmi following Dr Hou advice" edit:
# - Re-derive rows 2-17 of Sheet1 for the Spp1-Itgav LR pair table, now
#   cross-joining ALL 4 clusters (ECs, FAPs, M2, sCs) as both sending and
#   target cluster (was only FAPs/M2/sCs x ECs/FAPs/M2/sCs before, 12 rows).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array is: row number, then values for columns A..T in order.
$rows = @(
  ,@(2, "ECs", "Spp1", "Itgav", "ECs", 2, 0.6666666666666666, 12.486902, 37.460706, 0.01504353194025314, 0.01504353194025314, 3, 1, 13.441269, 40.323807, 0.08973082133481231, 0.08973082133481232, 167.839808758638, 1510.558278827742, 0.001349868476775397, 0.001349868476775397)
  ,@(3, "ECs", "Spp1", "Itgav", "FAPs", 2, 0.6666666666666666, 12.486902, 37.460706, 0.01504353194025314, 0.01504353194025314, 3, 1, 54.711535, 164.134605, 0.3652416280068742, 0.3652416280068742, 683.17757581457, 6148.59818233113, 0.005494524096831466, 0.005494524096831468)
  ,@(4, "ECs", "Spp1", "Itgav", "M2", 2, 0.6666666666666666, 12.486902, 37.460706, 0.01504353194025314, 0.01504353194025314, 3, 1, 63.67711, 191.03133, 0.4250937452800914, 0.4250937452800915, 795.12983221322, 7156.16848991898, 0.006394911334722887, 0.006394911334722889)
  ,@(5, "ECs", "Spp1", "Itgav", "sCs", 2, 0.6666666666666666, 12.486902, 37.460706, 0.01504353194025314, 0.01504353194025314, 3, 1, 17.96553866666667, 53.896616, 0.119933805378222, 0.119933805378222, 224.3339207078774, 2019.005286370896, 0.001804228031923386, 0.001804228031923387)
  ,@(6, "FAPs", "Spp1", "Itgav", "ECs", 3, 1, 46.08534733333334, 138.256042, 0.05552108878460485, 0.05552108878460485, 3, 1, 13.441269, 40.323807, 0.08973082133481231, 0.08973082133481232, 619.445550465766, 5575.009954191894, 0.004981952898045629, 0.00498195289804563)
  ,@(7, "FAPs", "Spp1", "Itgav", "FAPs", 3, 1, 46.08534733333334, 138.256042, 0.05552108878460485, 0.05552108878460485, 3, 1, 54.711535, 164.134605, 0.3652416280068742, 0.3652416280068742, 2521.400093614824, 22692.60084253341, 0.02027861285640328, 0.02027861285640328)
  ,@(8, "FAPs", "Spp1", "Itgav", "M2", 3, 1, 46.08534733333334, 138.256042, 0.05552108878460485, 0.05552108878460485, 3, 1, 63.67711, 191.03133, 0.4250937452800914, 0.4250937452800915, 2934.581731532874, 26411.23558379586, 0.02360166757347616, 0.02360166757347616)
  ,@(9, "FAPs", "Spp1", "Itgav", "sCs", 3, 1, 46.08534733333334, 138.256042, 0.05552108878460485, 0.05552108878460485, 3, 1, 17.96553866666667, 53.896616, 0.119933805378222, 0.119933805378222, 827.9480894837636, 7451.532805353872, 0.006658855456679783, 0.006658855456679783)
  ,@(10, "M2", "Spp1", "Itgav", "ECs", 3, 1, 666.4749603333333, 1999.424881, 0.8029323328679479, 0.8029323328679479, 3, 1, 13.441269, 40.323807, 0.08973082133481231, 0.08973082133481232, 8958.269223604662, 80624.42301244197, 0.07204777770451787, 0.07204777770451788)
  ,@(11, "M2", "Spp1", "Itgav", "FAPs", 3, 1, 666.4749603333333, 1999.424881, 0.8029323328679479, 0.8029323328679479, 3, 1, 54.711535, 164.134605, 0.3652416280068742, 0.3652416280068742, 36463.86811890078, 328174.813070107, 0.2932643124360467, 0.2932643124360467)
  ,@(12, "M2", "Spp1", "Itgav", "M2", 3, 1, 666.4749603333333, 1999.424881, 0.8029323328679479, 0.8029323328679479, 3, 1, 63.67711, 191.03133, 0.4250937452800914, 0.4250937452800915, 42439.1993613913, 381952.7942525217, 0.341321512585317, 0.3413215125853171)
  ,@(13, "M2", "Spp1", "Itgav", "sCs", 3, 1, 666.4749603333333, 1999.424881, 0.8029323328679479, 0.8029323328679479, 3, 1, 17.96553866666667, 53.896616, 0.119933805378222, 0.119933805378222, 11973.58167023363, 107762.2350321027, 0.09629873014206622, 0.09629873014206623)
  ,@(14, "sCs", "Spp1", "Itgav", "ECs", 3, 1, 105.0040076666667, 315.012023, 0.1265030464071941, 0.1265030464071941, 3, 1, 13.441269, 40.323807, 0.08973082133481231, 0.08973082133481232, 1411.387113125729, 12702.48401813156, 0.01135122225547341, 0.01135122225547341)
  ,@(15, "sCs", "Spp1", "Itgav", "FAPs", 3, 1, 105.0040076666667, 315.012023, 0.1265030464071941, 0.1265030464071941, 3, 1, 54.711535, 164.134605, 0.3652416280068742, 0.3652416280068742, 5744.930440595102, 51704.37396535591, 0.04620417861759274, 0.04620417861759275)
  ,@(16, "sCs", "Spp1", "Itgav", "M2", 3, 1, 105.0040076666667, 315.012023, 0.1265030464071941, 0.1265030464071941, 3, 1, 63.67711, 191.03133, 0.4250937452800914, 0.4250937452800915, 6686.351746631177, 60177.16571968059, 0.05377565378657537, 0.05377565378657538)
  ,@(17, "sCs", "Spp1", "Itgav", "sCs", 3, 1, 105.0040076666667, 315.012023, 0.1265030464071941, 0.1265030464071941, 3, 1, 17.96553866666667, 53.896616, 0.119933805378222, 0.119933805378222, 1886.453559890463, 16978.08203901417, 0.01517199174755261, 0.01517199174755261)
)

foreach ($entry in $rows) {
    $rownum = $entry[0]
    for ($i = 1; $i -lt $entry.Count; $i++) {
        $ws.Cells.Item($rownum, $i).Value = $entry[$i]
    }
}
